# regen sval data to filter save games
# Updates numeric cells B2:E7 and the sum column G2:G7 on the active sheet
# to the regenerated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.1190320826869504, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 1.069511820747223)
    3 = @(0.6606524410359556, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 2.960089034096801)
    4 = @(3.286832544864788, 0.04071648406533734, 0.7527432677738641, 10.19245300693656, 14.27274530364055)
    5 = @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 3.56341032713086)
    6 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    7 = @(0.04271373187048222, 0.04071648406533734, 22.3905356188092, 0.4942365360607697, 22.96820237080578)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G - sum
}
